$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.188.77'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.81'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("E4").Value = '  +0.47%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.27'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  +3.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  +0.99%  '

$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.04'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = '  +6.32%  '

$ws.Range("E9").Value = '  +2.07%  '

$ws.Range("E10").Value = '  +1.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0989'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.128.76'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").Value = '  +1.58%  '

$ws.Range("B13").Value = 'Chainlink'

$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.47'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("B14").Value = 'WrappedEther'

$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.46'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").Value = '  +1.63%  '

$ws.Range("E15").Value = '  +1.56%  '

$ws.Range("E16").Value = '  +1.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.151.18'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  +0.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.82'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0795'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").Value = '  +1.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.41'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  +0.48%  '

$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.72'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").Value = '  +1.37%  '

$ws.Range("E23").Value = '  +0.46%  '

$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.29'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").Value = '  -2.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.93'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = '  +27.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.98'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  +3.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.64'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").Value = '  +1.92%  '

$ws.Range("E29").Value = '  +0.35%  '

$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("E31").Value = '  +1.29%  '

$ws.Range("E32").Value = '  +2.58%  '

$ws.Range("E33").Value = '  +27.61%  '

$ws.Range("E34").Value = '  +2.28%  '

$ws.Range("E35").Value = '  +18.05%  '

$ws.Range("E36").Value = '  +9.40%  '

$ws.Range("E37").Value = '  +6.56%  '

$ws.Range("E38").Value = '  +5.05%  '

$ws.Range("E39").Value = '  +3.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '89.42'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").Value = '  -2.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.339.76'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("B42").Value = 'Kaspa'

$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0592'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  +13.41%  '

$ws.Range("B43").Value = 'InjectiveProtocol'

$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.85'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  +2.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.31'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").Value = '  +2.96%  '

$ws.Range("E45").Value = '  -0.10%  '

$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("B47").Value = 'Gas'

$ws.Range("C47").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.24'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").Value = '  +44.82%  '

$ws.Range("B48").Value = 'FraxShare'

$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.58'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").Value = '  +5.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.041.78'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  +1.55%  '

$ws.Range("E50").Value = '  +1.27%  '

$ws.Range("E51").Value = '  +0.62%  '
